# clean source text before match
# Row 26 matched the wrong CBDB person (李心傳, id 10831) because the
# "contents" text still contained the raw book title fragment that was
# being picked up by the writing-match step. Cleaning that source text
# before matching produces the correct person (鄭涇, id 51241) and moves
# the matched phrase from writing_match into source_match instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

# input_id / cbdb_id / person_name now point at the correct CBDB record
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "51241"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 51241
$ws.Cells.Item($row, 3).Value = "鄭涇"

# contents text cleaned up
$ws.Cells.Item($row, 11).Value = "見登科錄"

# the matched phrase moved from writing_match (T/U) to source_match (R/S)
$ws.Cells.Item($row, 18).Value = 1
$ws.Cells.Item($row, 19).Value = "登科錄"

$ws.Cells.Item($row, 20).Value = 0
$ws.Cells.Item($row, 21).Value = ""
